$wb = $excel.ActiveWorkbook

# --- Sheet: ip_address_list ---
$ws1 = $wb.Worksheets.Item("ip_address_list")

$ws1.Range("A1").Value = "kartoffelnsalat"
$ws1.Range("B1").Value = "192.168.000.j"
$ws1.Range("C1").Value = "255.255.255.0"
$ws1.Range("D1").Value = "kkgg"
$ws1.Range("E1").Value = 0

$ws1.Range("A2").Value = "einkaufenfh"
$ws1.Range("B2").Value = "192.168.000.000"
$ws1.Range("C2").Value = "255.255.255.0"
$ws1.Range("D2").Value = "gggg"
$ws1.Range("E2").Value = 1

$ws1.Range("A3").Value = "regenschrim"
$ws1.Range("B3").Value = "192.168.000.000"
$ws1.Range("C3").Value = "255.255.255.0"
$ws1.Range("D3").ClearContents()
$ws1.Range("E3").Value = 0

$ws1.Range("A4").Value = "bewolktt"
$ws1.Range("B4").Value = "192.168.000.000"
$ws1.Range("C4").Value = "255.255.255.0"
$ws1.Range("D4").Value = "du hast einen problem"
$ws1.Range("E4").Value = 1

$ws1.Range("A5").Value = "hggh"
$ws1.Range("B5").Value = "192.168.000.000"
$ws1.Range("C5").Value = "255.255.255.0"
$ws1.Range("D5").Value = "joo				"
$ws1.Range("E5").Value = 1

# --- Sheet: ip_adress_fav_list ---
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")

$ws2.Range("A1").Value = "einkaufenfh"
$ws2.Range("B1").Value = "192.168.000.000"
$ws2.Range("C1").Value = "255.255.255.0"
$ws2.Range("D1").Value = "gggg"
$ws2.Range("E1").Value = 1

$ws2.Range("A2").Value = "hggh"
$ws2.Range("B2").Value = "192.168.000.000"
$ws2.Range("C2").Value = "255.255.255.0"
$ws2.Range("D2").Value = "joo				"
$ws2.Range("E2").Value = 1

$ws2.Range("A3").Value = "bewolktt"
$ws2.Range("B3").Value = "192.168.000.000"
$ws2.Range("C3").Value = "255.255.255.0"
$ws2.Range("D3").Value = "du hast einen problem"
$ws2.Range("E3").Value = 1

# --- Sheet: Settings ---
$ws4 = $wb.Worksheets.Item("Settings")
$ws4.Range("B4").Value = 0
